$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 10 and 11 (the "0xAX" / "0xBX" placeholder rows), which
# shifts all subsequent rows up by two and removes the now-unused
# shared-string entries for "0xAX" and "0xBX".
$ws.Rows("10:11").Delete()

# Update the active selection to B6, matching the saved view state.
$ws.Activate()
$ws.Range("B6").Select()
